# Mini Fixes for staging
# billing/Working/Prod/rates/rates.xlsx
#
# The rate for "Navel" / "All" (cell B2 on the "rates" sheet) changes
# from 140 to 400.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 400
